$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all the words (rows 1-378) in column C as "Learned" instead of
# "No learned". This also selects the range, matching the committed
# selection state (activeCell C1, sqref C1:C378).
$rng = $ws.Range("C1:C378")
$rng.Select()
$rng.Value = "Learned"
